# week03_report.xlsx edit: mark PPP014 (row 11) as Completed and refresh the
# "Generated:" timestamp stamp, matching the authoring tool's output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (student PPP014 / Srinivasan SR) flips from "Pending" to "Completed":
#  - clear the Pending Task cell (D11) entirely
#  - change the Completion Status cell (E11) text to "Completed"
#  - apply the same "Completed" look (bold white on green) used by the other
#    completed rows to B11 and E11, by copying the formatting from row 10
#    (already in the "Completed" state) instead of re-building it from
#    scratch so we reuse the existing style entries.

$ws.Range("D11").ClearContents()
$ws.Range("E11").Value = "Completed"

$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)

$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Refresh the "Generated: <date>" footer stamp.
$ws.Range("A30").Value = "Generated: 2023-09-14 10:47:40 AM"
